$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.030314852683073
$ws.Cells.Item(2, 4).Value = 1.040630672165149
$ws.Cells.Item(2, 5).Value = 1.030038938321472
$ws.Cells.Item(2, 6).Value = 1.049332779181627
$ws.Cells.Item(2, 9).Value = 1.038347399260381
$ws.Cells.Item(2, 10).Value = 1.035456778597277
$ws.Cells.Item(2, 11).Value = 1.043412474603123
$ws.Cells.Item(2, 12).Value = 1.032851106593954
$ws.Cells.Item(2, 13).Value = 1.052090130802789
$ws.Cells.Item(2, 14).Value = 1.01581936893066

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.031122909739938
$ws.Cells.Item(3, 4).Value = 1.041306028644359
$ws.Cells.Item(3, 5).Value = 1.030721595259665
$ws.Cells.Item(3, 6).Value = 1.050240714774217
$ws.Cells.Item(3, 9).Value = 1.038560865243187
$ws.Cells.Item(3, 10).Value = 1.035906959674534
$ws.Cells.Item(3, 11).Value = 1.043898649231607
$ws.Cells.Item(3, 12).Value = 1.033342372624647
$ws.Cells.Item(3, 13).Value = 1.05281004412582
$ws.Cells.Item(3, 14).Value = 1.015969251711364

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.031646311488617
$ws.Cells.Item(4, 4).Value = 1.041743444932579
$ws.Cells.Item(4, 5).Value = 1.031164169607592
$ws.Cells.Item(4, 6).Value = 1.050829169929479
$ws.Cells.Item(4, 9).Value = 1.038697887169163
$ws.Cells.Item(4, 10).Value = 1.036198123807216
$ws.Cells.Item(4, 11).Value = 1.044212971293788
$ws.Cells.Item(4, 12).Value = 1.033660425454589
$ws.Cells.Item(4, 13).Value = 1.05327620269127
$ws.Cells.Item(4, 14).Value = 1.01606616502911

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.031866475971954
$ws.Cells.Item(5, 4).Value = 1.041927432434479
$ws.Cells.Item(5, 5).Value = 1.031350429477989
$ws.Cells.Item(5, 6).Value = 1.051076784450203
$ws.Cells.Item(5, 9).Value = 1.038755225874004
$ws.Cells.Item(5, 10).Value = 1.03632049612302
$ws.Cells.Item(5, 11).Value = 1.044345047402572
$ws.Cells.Item(5, 12).Value = 1.033794174434312
$ws.Cells.Item(5, 13).Value = 1.05347225226578
$ws.Cells.Item(5, 14).Value = 1.016106889982332

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.031903449921188
$ws.Cells.Item(6, 4).Value = 1.041958330422861
$ws.Cells.Item(6, 5).Value = 1.031381715100313
$ws.Cells.Item(6, 6).Value = 1.051118373331623
$ws.Cells.Item(6, 9).Value = 1.038764837713188
$ws.Cells.Item(6, 10).Value = 1.03634104099783
$ws.Cells.Item(6, 11).Value = 1.044367219718444
$ws.Cells.Item(6, 12).Value = 1.033816633758134
$ws.Cells.Item(6, 13).Value = 1.053505174290756
$ws.Cells.Item(6, 14).Value = 1.016113726847686

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.031649252842656
$ws.Cells.Item(7, 4).Value = 1.041745903000942
$ws.Cells.Item(7, 5).Value = 1.03116665763115
$ws.Cells.Item(7, 6).Value = 1.0508324776734
$ws.Cells.Item(7, 9).Value = 1.038698654375255
$ws.Cells.Item(7, 10).Value = 1.036199759083617
$ws.Cells.Item(7, 11).Value = 1.044214736357859
$ws.Cells.Item(7, 12).Value = 1.033662212461195
$ws.Cells.Item(7, 13).Value = 1.053278822016706
$ws.Cells.Item(7, 14).Value = 1.016066709266873

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.030587827626033
$ws.Cells.Item(8, 4).Value = 1.040858824954103
$ws.Cells.Item(8, 5).Value = 1.030269468422931
$ws.Cells.Item(8, 6).Value = 1.049639420337355
$ws.Cells.Item(8, 9).Value = 1.038419769364344
$ws.Cells.Item(8, 10).Value = 1.035608946234741
$ws.Cells.Item(8, 11).Value = 1.043576833748285
$ws.Cells.Item(8, 12).Value = 1.033017096070309
$ws.Cells.Item(8, 13).Value = 1.052333360281904
$ws.Cells.Item(8, 14).Value = 1.015870036926159

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.028721629434727
$ws.Cells.Item(9, 4).Value = 1.03929894059892
$ws.Cells.Item(9, 5).Value = 1.028695093142135
$ws.Cells.Item(9, 6).Value = 1.047544537358237
$ws.Cells.Item(9, 9).Value = 1.037919912523615
$ws.Cells.Item(9, 10).Value = 1.034566892596365
$ws.Cells.Item(9, 11).Value = 1.042450794230134
$ws.Cells.Item(9, 12).Value = 1.031881686064854
$ws.Cells.Item(9, 13).Value = 1.050669909628968
$ws.Cells.Item(9, 14).Value = 1.015522951866767

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.027480397904096
$ws.Cells.Item(10, 4).Value = 1.038261315537718
$ws.Cells.Item(10, 5).Value = 1.027650045270154
$ws.Cells.Item(10, 6).Value = 1.046153060512207
$ws.Cells.Item(10, 9).Value = 1.037581057941691
$ws.Cells.Item(10, 10).Value = 1.033871613054649
$ws.Cells.Item(10, 11).Value = 1.0416988541065
$ws.Cells.Item(10, 12).Value = 1.031125744948803
$ws.Cells.Item(10, 13).Value = 1.04956277133808
$ws.Cells.Item(10, 14).Value = 1.015291236242134

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.026943639489754
$ws.Cells.Item(11, 4).Value = 1.037812579740736
$ws.Cells.Item(11, 5).Value = 1.027198626440808
$ws.Cells.Item(11, 6).Value = 1.045551771689979
$ws.Cells.Item(11, 9).Value = 1.037433008119195
$ws.Cells.Item(11, 10).Value = 1.033570427679792
$ws.Cells.Item(11, 11).Value = 1.041372976742191
$ws.Cells.Item(11, 12).Value = 1.030798667786003
$ws.Cells.Item(11, 13).Value = 1.049083822711439
$ws.Cells.Item(11, 14).Value = 1.015190829121604

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.026744370785342
$ws.Cells.Item(12, 4).Value = 1.037645985625314
$ws.Cells.Item(12, 5).Value = 1.02703111555762
$ws.Cells.Item(12, 6).Value = 1.045328612915767
$ws.Cells.Item(12, 9).Value = 1.037377817772978
$ws.Cells.Item(12, 10).Value = 1.033458536500384
$ws.Cells.Item(12, 11).Value = 1.041251890308372
$ws.Cells.Item(12, 12).Value = 1.030677215725524
$ws.Cells.Item(12, 13).Value = 1.048905989070161
$ws.Cells.Item(12, 14).Value = 1.015153522963823

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.02678710971821
$ws.Cells.Item(13, 4).Value = 1.037681716692667
$ws.Cells.Item(13, 5).Value = 1.027067039661681
$ws.Cells.Item(13, 6).Value = 1.04537647275282
$ws.Cells.Item(13, 9).Value = 1.037389665238945
$ws.Cells.Item(13, 10).Value = 1.033482538311455
$ws.Cells.Item(13, 11).Value = 1.041277865596413
$ws.Cells.Item(13, 12).Value = 1.030703265820023
$ws.Cells.Item(13, 13).Value = 1.048944131831567
$ws.Cells.Item(13, 14).Value = 1.015161525728547

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.026927165678675
$ws.Cells.Item(14, 4).Value = 1.037798807247484
$ws.Cells.Item(14, 5).Value = 1.027184776538518
$ws.Cells.Item(14, 6).Value = 1.045533321497163
$ws.Cells.Item(14, 9).Value = 1.0374284501035
$ws.Cells.Item(14, 10).Value = 1.03356117907036
$ws.Cells.Item(14, 11).Value = 1.041362968529853
$ws.Cells.Item(14, 12).Value = 1.030788627718194
$ws.Cells.Item(14, 13).Value = 1.04906912150582
$ws.Cells.Item(14, 14).Value = 1.015187745594242

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.027013473041966
$ws.Cells.Item(15, 4).Value = 1.037870962097414
$ws.Cells.Item(15, 5).Value = 1.027257340188165
$ws.Cells.Item(15, 6).Value = 1.045629985977162
$ws.Cells.Item(15, 9).Value = 1.037452320518951
$ws.Cells.Item(15, 10).Value = 1.03360962995221
$ws.Cells.Item(15, 11).Value = 1.041415397853836
$ws.Cells.Item(15, 12).Value = 1.030841227204063
$ws.Cells.Item(15, 13).Value = 1.049146140989233
$ws.Cells.Item(15, 14).Value = 1.01520389914308

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.027516035744894
$ws.Cells.Item(16, 4).Value = 1.038291108663858
$ws.Cells.Item(16, 5).Value = 1.02768002765321
$ws.Cells.Item(16, 6).Value = 1.046192992144231
$ws.Cells.Item(16, 9).Value = 1.037590855703224
$ws.Cells.Item(16, 10).Value = 1.033891599212012
$ws.Cells.Item(16, 11).Value = 1.041720475680447
$ws.Cells.Item(16, 12).Value = 1.031147457372458
$ws.Cells.Item(16, 13).Value = 1.049594567194656
$ws.Cells.Item(16, 14).Value = 1.015297898437669

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.027831469434273
$ws.Cells.Item(17, 4).Value = 1.038554807427419
$ws.Cells.Item(17, 5).Value = 1.027945462472103
$ws.Cells.Item(17, 6).Value = 1.046546481639606
$ws.Cells.Item(17, 9).Value = 1.037677401375009
$ws.Cells.Item(17, 10).Value = 1.034068438543158
$ws.Cells.Item(17, 11).Value = 1.041911768303664
$ws.Cells.Item(17, 12).Value = 1.031339615441482
$ws.Cells.Item(17, 13).Value = 1.049875974639393
$ws.Cells.Item(17, 14).Value = 1.015356842554274

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.028015524290103
$ws.Cells.Item(18, 4).Value = 1.038708672626748
$ws.Cells.Item(18, 5).Value = 1.028100391523405
$ws.Cells.Item(18, 6).Value = 1.046752784615545
$ws.Cells.Item(18, 9).Value = 1.037727754243881
$ws.Cells.Item(18, 10).Value = 1.034171573815199
$ws.Cells.Item(18, 11).Value = 1.042023318723517
$ws.Cells.Item(18, 12).Value = 1.031451722054618
$ws.Cells.Item(18, 13).Value = 1.050040158105295
$ws.Cells.Item(18, 14).Value = 1.015391216613505

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.028078293676188
$ws.Cells.Item(19, 4).Value = 1.038761145817201
$ws.Cells.Item(19, 5).Value = 1.028153236149924
$ws.Cells.Item(19, 6).Value = 1.046823148646663
$ws.Cells.Item(19, 9).Value = 1.037744901576812
$ws.Cells.Item(19, 10).Value = 1.034206738185157
$ws.Cells.Item(19, 11).Value = 1.042061349857235
$ws.Cells.Item(19, 12).Value = 1.031489951552879
$ws.Cells.Item(19, 13).Value = 1.05009614771211
$ws.Cells.Item(19, 14).Value = 1.015402936070645

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.02779761936799
$ws.Cells.Item(20, 4).Value = 1.038526509429458
$ws.Cells.Item(20, 5).Value = 1.027916972911534
$ws.Cells.Item(20, 6).Value = 1.04650854328793
$ws.Cells.Item(20, 9).Value = 1.037668129052255
$ws.Cells.Item(20, 10).Value = 1.034049466601122
$ws.Cells.Item(20, 11).Value = 1.041891247232814
$ws.Cells.Item(20, 12).Value = 1.031318996200586
$ws.Cells.Item(20, 13).Value = 1.049845777796743
$ws.Cells.Item(20, 14).Value = 1.015350519134078

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.026885919714707
$ws.Cells.Item(21, 4).Value = 1.037764324599658
$ws.Cells.Item(21, 5).Value = 1.027150101360961
$ws.Cells.Item(21, 6).Value = 1.045487128283144
$ws.Cells.Item(21, 9).Value = 1.037417034388312
$ws.Cells.Item(21, 10).Value = 1.033538021800795
$ws.Cells.Item(21, 11).Value = 1.041337908960047
$ws.Cells.Item(21, 12).Value = 1.030763489680412
$ws.Cells.Item(21, 13).Value = 1.049032313238319
$ws.Cells.Item(21, 14).Value = 1.015180024783118

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.026313319404869
$ws.Cells.Item(22, 4).Value = 1.037285609169547
$ws.Cells.Item(22, 5).Value = 1.026668901093365
$ws.Cells.Item(22, 6).Value = 1.044846005066271
$ws.Cells.Item(22, 9).Value = 1.037258015745963
$ws.Cells.Item(22, 10).Value = 1.03321635538903
$ws.Cells.Item(22, 11).Value = 1.040989766267667
$ws.Cells.Item(22, 12).Value = 1.030414446767044
$ws.Cells.Item(22, 13).Value = 1.048521256999683
$ws.Cells.Item(22, 14).Value = 1.015072767788021

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.026616806206283
$ws.Cells.Item(23, 4).Value = 1.037539337241802
$ws.Cells.Item(23, 5).Value = 1.026923902630376
$ws.Cells.Item(23, 6).Value = 1.045185773519665
$ws.Cells.Item(23, 9).Value = 1.037342422846476
$ws.Cells.Item(23, 10).Value = 1.033386885984056
$ws.Cells.Item(23, 11).Value = 1.041174345307619
$ws.Cells.Item(23, 12).Value = 1.030599459138372
$ws.Cells.Item(23, 13).Value = 1.048792138995352
$ws.Cells.Item(23, 14).Value = 1.015129632330508

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.027812914546834
$ws.Cells.Item(24, 4).Value = 1.038539295909509
$ws.Cells.Item(24, 5).Value = 1.027929845792169
$ws.Cells.Item(24, 6).Value = 1.046525685629995
$ws.Cells.Item(24, 9).Value = 1.037672319210626
$ws.Cells.Item(24, 10).Value = 1.034058039242682
$ws.Cells.Item(24, 11).Value = 1.041900519906023
$ws.Cells.Item(24, 12).Value = 1.031328313073823
$ws.Cells.Item(24, 13).Value = 1.049859422316649
$ws.Cells.Item(24, 14).Value = 1.015353376437304

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.029203582319654
$ws.Cells.Item(25, 4).Value = 1.039701810813525
$ws.Cells.Item(25, 5).Value = 1.029101315281616
$ws.Cells.Item(25, 6).Value = 1.048085222030495
$ws.Cells.Item(25, 9).Value = 1.038050131023289
$ws.Cells.Item(25, 10).Value = 1.034836395453075
$ws.Cells.Item(25, 11).Value = 1.042742128092722
$ws.Cells.Item(25, 12).Value = 1.03217504651884
$ws.Cells.Item(25, 13).Value = 1.05109963636027
$ws.Cells.Item(25, 14).Value = 1.015612740979615

